$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 169
$ws.Range("F3").Value = 377
$ws.Range("F4").Value = 1068
$ws.Range("F5").Value = 23
$ws.Range("F9").Value = 293
$ws.Range("F10").Value = 393
$ws.Range("F13").Value = 329
$ws.Range("F15").Value = 299
$ws.Range("F16").Value = 412
$ws.Range("F17").Value = 5368
$ws.Range("F18").Value = 82
$ws.Range("F19").Value = 1492
$ws.Range("F20").Value = 335
$ws.Range("F21").Value = 4437
$ws.Range("F22").Value = 107
$ws.Range("F23").Value = 76
$ws.Range("F24").Value = 1419
$ws.Range("F27").Value = 621
$ws.Range("F29").Value = 3763

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 3748
$ws.Range("F3").Value = 3748
$ws.Range("F6").Value = 93

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9339
$ws.Range("F3").Value = 573
$ws.Range("F4").Value = 2090

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9339
$ws.Range("F3").Value = 573
$ws.Range("F4").Value = 2090
$ws.Range("F5").Value = 3748
$ws.Range("F6").Value = 3748
$ws.Range("F7").Value = 169
$ws.Range("F8").Value = 378
$ws.Range("F9").Value = 1068
$ws.Range("F10").Value = 23
$ws.Range("F14").Value = 293
$ws.Range("F15").Value = 393
$ws.Range("F18").Value = 329
$ws.Range("F22").Value = 93
$ws.Range("F24").Value = 299
$ws.Range("F26").Value = 412
$ws.Range("F27").Value = 5368
$ws.Range("F28").Value = 82
$ws.Range("F29").Value = 1492
$ws.Range("F32").Value = 335
$ws.Range("F34").Value = 4437
$ws.Range("F35").Value = 107
$ws.Range("F36").Value = 76
$ws.Range("F37").Value = 1419
$ws.Range("F40").Value = 621
$ws.Range("F47").Value = 3763
